# Applies the scheduled-runner profit/price recalculation update to the Typhon sheets.
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) gets updated average-price /
# leve-profit figures in columns H:N for the affected rows. A few cells (previously
# holding a stray LeveProfit value) are cleared outright because that column no longer
# applies to those rows after the recalculation.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 73.72
$ws.Range("I15").Value = 73.72
$ws.Range("K15").Value = 221.16
$ws.Range("M15").Value = -52.16
$ws.Range("H21").Value = 35000
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = ""
$ws.Range("H23").Value = 35000
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = ""
$ws.Range("H129").Value = 1445.5333
$ws.Range("J129").Value = 1445.5333
$ws.Range("L129").Value = 4336.5999
$ws.Range("N129").Value = -14336.5999
$ws.Range("H135").Value = 13516782
$ws.Range("I135").Value = 632.5806
$ws.Range("J135").Value = 83350220
$ws.Range("K135").Value = 5693.2254
$ws.Range("L135").Value = 750151980
$ws.Range("M135").Value = -3158.2254
$ws.Range("N135").Value = -750157050
$ws.Range("H138").Value = 121688.12
$ws.Range("I138").Value = 2213.2727
$ws.Range("J138").Value = 139447.89
$ws.Range("K138").Value = 6639.8181
$ws.Range("L138").Value = 418343.67
$ws.Range("M138").Value = -1499.8181
$ws.Range("N138").Value = -428623.67

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2770.1428
$ws.Range("J61").Value = 3579.8
$ws.Range("L61").Value = 3579.8
$ws.Range("N61").Value = -4003.8
$ws.Range("H102").Value = 1482.0555
$ws.Range("I102").Value = 1378.4667
$ws.Range("K102").Value = 1378.4667
$ws.Range("M102").Value = 243.5333000000001
$ws.Range("H132").Value = 12863.319
$ws.Range("I132").Value = 1841.1621
$ws.Range("K132").Value = 5523.4863
$ws.Range("M132").Value = -2993.4863
$ws.Range("H136").Value = 2770.1428
$ws.Range("J136").Value = 3579.8
$ws.Range("L136").Value = 10739.4
$ws.Range("N136").Value = -15839.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1165.1538
$ws.Range("I94").Value = 882.8333
$ws.Range("J94").Value = 1407.1428
$ws.Range("K94").Value = 882.8333
$ws.Range("L94").Value = 1407.1428
$ws.Range("M94").Value = -431.8333
$ws.Range("N94").Value = -2309.1428
$ws.Range("H105").Value = 1138404.2
$ws.Range("I105").Value = 1644.375
$ws.Range("J105").Value = 1787981.2
$ws.Range("K105").Value = 1644.375
$ws.Range("L105").Value = 1787981.2
$ws.Range("M105").Value = 102.625
$ws.Range("N105").Value = -1791475.2
$ws.Range("H134").Value = 3826.457
$ws.Range("I134").Value = 3792.5293
$ws.Range("K134").Value = 11377.5879
$ws.Range("M134").Value = -8842.5879

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 17860682
$ws.Range("I99").Value = 2971.2144
$ws.Range("K99").Value = 2971.2144
$ws.Range("M99").Value = -1473.2144
$ws.Range("H122").Value = 911.125
$ws.Range("I122").Value = 911.125
$ws.Range("K122").Value = 2733.375
$ws.Range("M122").Value = -283.375
$ws.Range("H126").Value = 17860682
$ws.Range("I126").Value = 2971.2144
$ws.Range("K126").Value = 8913.643199999999
$ws.Range("M126").Value = -6443.643199999999
$ws.Range("H132").Value = 3494.2
$ws.Range("I132").Value = 2364.8
$ws.Range("J132").Value = 6882.4
$ws.Range("K132").Value = 7094.400000000001
$ws.Range("L132").Value = 20647.2
$ws.Range("M132").Value = -4564.400000000001
$ws.Range("N132").Value = -25707.2
$ws.Range("H134").Value = 1463.6154
$ws.Range("I134").Value = 1131.2
$ws.Range("J134").Value = 2571.6667
$ws.Range("K134").Value = 3393.6
$ws.Range("L134").Value = 7715.000100000001
$ws.Range("M134").Value = -858.6000000000004
$ws.Range("N134").Value = -12785.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 7507.5
$ws.Range("J81").Value = 7507.5
$ws.Range("L81").Value = 22522.5
$ws.Range("N81").Value = -24768.5
$ws.Range("H84").Value = 7507.5
$ws.Range("J84").Value = 7507.5
$ws.Range("L84").Value = 67567.5
$ws.Range("N84").Value = -78799.5
$ws.Range("H92").Value = 15625594
$ws.Range("I92").Value = 25000250
$ws.Range("K92").Value = 75000750
$ws.Range("M92").Value = -74999502
$ws.Range("H131").Value = 750.74
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 750.74
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2252.22
$ws.Range("M131").Value = ""
$ws.Range("N131").Value = -12332.22
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 148277.72
$ws.Range("I43").Value = 5837.5386
$ws.Range("K43").Value = 5837.5386
$ws.Range("M43").Value = -5686.5386
$ws.Range("H46").Value = 34800
$ws.Range("J46").Value = 34800
$ws.Range("L46").Value = 34800
$ws.Range("N46").Value = -35112
$ws.Range("H57").Value = 28440
$ws.Range("J57").Value = 28440
$ws.Range("L57").Value = 28440
$ws.Range("N57").Value = -30080
$ws.Range("H102").Value = 2354.4
$ws.Range("I102").Value = 2145.0476
$ws.Range("K102").Value = 2145.0476
$ws.Range("M102").Value = -523.0475999999999
$ws.Range("H126").Value = 5204.125
$ws.Range("J126").Value = 4449.9165
$ws.Range("L126").Value = 13349.7495
$ws.Range("N126").Value = -18289.7495
$ws.Range("H132").Value = 13039.24
$ws.Range("I132").Value = 3864.36
$ws.Range("J132").Value = 22214.12
$ws.Range("K132").Value = 11593.08
$ws.Range("L132").Value = 66642.36
$ws.Range("M132").Value = -9063.08
$ws.Range("N132").Value = -71702.36

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3139.8
$ws.Range("I22").Value = 10001
$ws.Range("J22").Value = 1424.5
$ws.Range("K22").Value = 10001
$ws.Range("L22").Value = 1424.5
$ws.Range("M22").Value = -9706
$ws.Range("N22").Value = -2014.5
$ws.Range("H27").Value = 3139.8
$ws.Range("I27").Value = 10001
$ws.Range("J27").Value = 1424.5
$ws.Range("K27").Value = 10001
$ws.Range("L27").Value = 1424.5
$ws.Range("M27").Value = -9894
$ws.Range("N27").Value = -1638.5
$ws.Range("H68").Value = 2199.5454
$ws.Range("I68").Value = 1425
$ws.Range("J68").Value = 2642.1428
$ws.Range("K68").Value = 1425
$ws.Range("L68").Value = 2642.1428
$ws.Range("M68").Value = -676
$ws.Range("N68").Value = -4140.1428
$ws.Range("H71").Value = 2199.5454
$ws.Range("I71").Value = 1425
$ws.Range("J71").Value = 2642.1428
$ws.Range("K71").Value = 7125
$ws.Range("L71").Value = 13210.714
$ws.Range("M71").Value = -3381
$ws.Range("N71").Value = -20698.714
$ws.Range("H122").Value = 1311521.9
$ws.Range("I122").Value = 1512294.5
$ws.Range("K122").Value = 4536883.5
$ws.Range("M122").Value = -4534433.5
$ws.Range("H132").Value = 253529.17
$ws.Range("I132").Value = 378392.12
$ws.Range("K132").Value = 1135176.36
$ws.Range("M132").Value = -1132646.36

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1023.25
$ws.Range("I122").Value = 1026.5714
$ws.Range("K122").Value = 3079.7142
$ws.Range("M122").Value = -629.7142000000003
$ws.Range("H132").Value = 914.93335
$ws.Range("I132").Value = 447.70834
$ws.Range("J132").Value = 2783.8333
$ws.Range("K132").Value = 1343.12502
$ws.Range("L132").Value = 8351.499899999999
$ws.Range("M132").Value = 1186.87498
$ws.Range("N132").Value = -13411.4999
